# "more work on report"
#
# - Mark backlog rows 11 & 12 ("Terminé") by stamping their "Date fait"
#   (column H) with today's date, copying the date number-format from the
#   already-completed row 10 so the style matches (fill/border + date
#   format) instead of picking up a brand new style.
# - Add a new backlog line (row 57): Contenu / Données / PIM / "Faire
#   l'analyse par type de produit" / 2 pts.
# - Re-point the UI at the backlog sheet (it was left on Burndown) and
#   land the selection on H14, Burndown's view scrolled back up a row.
# - Everything else (G11/G12 flipping to "Terminé", the E1 total, the
#   Burndown SUMPRODUCT/shared formulas for the "today" column, …) is
#   purely formula-driven off TODAY()/H11/H12/E57 and falls out of the
#   automatic recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("backlog")
$bd = $wb.Worksheets.Item("Burndown")

# --- Close out rows 11 & 12 ------------------------------------------------
# Copy H10's format (date number format + fill/border) onto H11/H12 so the
# cell style matches the rest of the "Date fait" column, then stamp today.
$ws.Range("H10").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("H12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H11").Value = 43955
$ws.Range("H12").Value = 43955

# --- New backlog entry (row 57) --------------------------------------------
$ws.Range("B57").Value = "Contenu"
$ws.Range("C57").Value = "Données / PIM"
$ws.Range("D57").Value = "Faire l'analyse par type de produit"
$ws.Range("E57").Value = 2

# --- Recalculate so G11/G12/G57, E1, and the Burndown sheet pick up the
#     new inputs ------------------------------------------------------------
$excel.Calculate()

# --- UI state: backlog becomes the active/selected tab, cursor on H14 ------
$bd.Activate()
$bd.Range("A3").Select()
$ws.Activate()
$ws.Range("H14").Select()
